# Appends 12 new daily COVID-19 data rows (2022/10/21 - 2022/11/01) to Sheet1
# as rows 937-948, mirroring the structure/columns of the preceding rows
# (A=Data, B=Dia desde contagio 1, C=Casos confirmados, D=Obitos,
# F=Taxa morte contaminados, G=Curados, H=Casos negativos, I=Testes realizados,
# J=Novos Casos, K=Novos obitos, L=Novos testes realizados, M=Suspeitos,
# O=Suspeitos ativos, P=Leitos clinicos ocupados, Q=Leitos UTI ocupados - COVID-19,
# T=Semana de pandemia).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 937: 2022/10/21
$cell = $ws.Range("A937")
$cell.Formula = "=TEXT(""2022/10/21"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B937").Value = 936
$ws.Range("C937").Value = 39445
$ws.Range("D937").Value = 690
$ws.Range("F937").Value = 0.01749271137
$ws.Range("G937").Value = 38683
$ws.Range("H937").Value = 27040
$ws.Range("I937").Value = 66485
$ws.Range("J937").Value = 0
$ws.Range("K937").Value = 0
$ws.Range("L937").Value = 0
$ws.Range("M937").Value = 0
$ws.Range("O937").Value = 0
$ws.Range("P937").Value = 1
$ws.Range("Q937").Value = 2
$ws.Range("T937").Value = 134

# Row 938: 2022/10/22
$cell = $ws.Range("A938")
$cell.Formula = "=TEXT(""2022/10/22"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B938").Value = 937
$ws.Range("C938").Value = 39445
$ws.Range("D938").Value = 690
$ws.Range("F938").Value = 0.01749271137
$ws.Range("G938").Value = 38683
$ws.Range("H938").Value = 27040
$ws.Range("I938").Value = 66485
$ws.Range("J938").Value = 0
$ws.Range("K938").Value = 0
$ws.Range("L938").Value = 0
$ws.Range("M938").Value = 0
$ws.Range("O938").Value = 0
$ws.Range("P938").Value = 1
$ws.Range("Q938").Value = 2
$ws.Range("T938").Value = 134

# Row 939: 2022/10/23
$cell = $ws.Range("A939")
$cell.Formula = "=TEXT(""2022/10/23"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B939").Value = 938
$ws.Range("C939").Value = 39445
$ws.Range("D939").Value = 690
$ws.Range("F939").Value = 0.01749271137
$ws.Range("G939").Value = 38683
$ws.Range("H939").Value = 27040
$ws.Range("I939").Value = 66485
$ws.Range("J939").Value = 0
$ws.Range("K939").Value = 0
$ws.Range("L939").Value = 0
$ws.Range("M939").Value = 0
$ws.Range("O939").Value = 0
$ws.Range("P939").Value = 1
$ws.Range("Q939").Value = 2
$ws.Range("T939").Value = 135

# Row 940: 2022/10/24
$cell = $ws.Range("A940")
$cell.Formula = "=TEXT(""2022/10/24"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B940").Value = 939
$ws.Range("C940").Value = 39445
$ws.Range("D940").Value = 690
$ws.Range("F940").Value = 0.01749271137
$ws.Range("G940").Value = 38683
$ws.Range("H940").Value = 27040
$ws.Range("I940").Value = 66485
$ws.Range("J940").Value = 0
$ws.Range("K940").Value = 0
$ws.Range("L940").Value = 0
$ws.Range("M940").Value = 0
$ws.Range("O940").Value = 0
$ws.Range("P940").Value = 1
$ws.Range("Q940").Value = 2
$ws.Range("T940").Value = 135

# Row 941: 2022/10/25
$cell = $ws.Range("A941")
$cell.Formula = "=TEXT(""2022/10/25"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B941").Value = 940
$ws.Range("C941").Value = 39447
$ws.Range("D941").Value = 690
$ws.Range("F941").Value = 0.01749182447
$ws.Range("G941").Value = 38684
$ws.Range("H941").Value = 27040
$ws.Range("I941").Value = 66487
$ws.Range("J941").Value = 2
$ws.Range("K941").Value = 0
$ws.Range("L941").Value = 2
$ws.Range("M941").Value = 0
$ws.Range("O941").Value = 0
$ws.Range("P941").Value = 0
$ws.Range("Q941").Value = 2
$ws.Range("T941").Value = 135

# Row 942: 2022/10/26
$cell = $ws.Range("A942")
$cell.Formula = "=TEXT(""2022/10/26"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B942").Value = 941
$ws.Range("C942").Value = 39447
$ws.Range("D942").Value = 690
$ws.Range("F942").Value = 0.01749182447
$ws.Range("G942").Value = 38684
$ws.Range("H942").Value = 27040
$ws.Range("I942").Value = 66487
$ws.Range("J942").Value = 2
$ws.Range("K942").Value = 0
$ws.Range("L942").Value = 2
$ws.Range("M942").Value = 0
$ws.Range("O942").Value = 0
$ws.Range("P942").Value = 0
$ws.Range("Q942").Value = 2
$ws.Range("T942").Value = 135

# Row 943: 2022/10/27
$cell = $ws.Range("A943")
$cell.Formula = "=TEXT(""2022/10/27"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B943").Value = 942
$ws.Range("C943").Value = 39447
$ws.Range("D943").Value = 690
$ws.Range("F943").Value = 0.01749182447
$ws.Range("G943").Value = 38684
$ws.Range("H943").Value = 27040
$ws.Range("I943").Value = 66487
$ws.Range("J943").Value = 0
$ws.Range("K943").Value = 0
$ws.Range("L943").Value = 0
$ws.Range("M943").Value = 2
$ws.Range("O943").Value = 2
$ws.Range("P943").Value = 0
$ws.Range("Q943").Value = 2
$ws.Range("T943").Value = 135

# Row 944: 2022/10/28
$cell = $ws.Range("A944")
$cell.Formula = "=TEXT(""2022/10/28"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B944").Value = 943
$ws.Range("C944").Value = 39447
$ws.Range("D944").Value = 690
$ws.Range("F944").Value = 0.01749182447
$ws.Range("G944").Value = 38684
$ws.Range("H944").Value = 27040
$ws.Range("I944").Value = 66487
$ws.Range("J944").Value = 0
$ws.Range("K944").Value = 0
$ws.Range("L944").Value = 0
$ws.Range("M944").Value = 2
$ws.Range("O944").Value = 2
$ws.Range("P944").Value = 0
$ws.Range("Q944").Value = 2
$ws.Range("T944").Value = 135

# Row 945: 2022/10/29
$cell = $ws.Range("A945")
$cell.Formula = "=TEXT(""2022/10/29"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B945").Value = 944
$ws.Range("C945").Value = 39447
$ws.Range("D945").Value = 690
$ws.Range("F945").Value = 0.01749182447
$ws.Range("G945").Value = 38684
$ws.Range("H945").Value = 27040
$ws.Range("I945").Value = 66487
$ws.Range("J945").Value = 0
$ws.Range("K945").Value = 0
$ws.Range("L945").Value = 0
$ws.Range("M945").Value = 2
$ws.Range("O945").Value = 2
$ws.Range("P945").Value = 0
$ws.Range("Q945").Value = 2
$ws.Range("T945").Value = 135

# Row 946: 2022/10/30
$cell = $ws.Range("A946")
$cell.Formula = "=TEXT(""2022/10/30"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B946").Value = 945
$ws.Range("C946").Value = 39447
$ws.Range("D946").Value = 690
$ws.Range("F946").Value = 0.01749182447
$ws.Range("G946").Value = 38684
$ws.Range("H946").Value = 27040
$ws.Range("I946").Value = 66487
$ws.Range("J946").Value = 0
$ws.Range("K946").Value = 0
$ws.Range("L946").Value = 0
$ws.Range("M946").Value = 2
$ws.Range("O946").Value = 2
$ws.Range("P946").Value = 0
$ws.Range("Q946").Value = 2
$ws.Range("T946").Value = 136

# Row 947: 2022/10/31
$cell = $ws.Range("A947")
$cell.Formula = "=TEXT(""2022/10/31"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B947").Value = 946
$ws.Range("C947").Value = 39447
$ws.Range("D947").Value = 690
$ws.Range("F947").Value = 0.01749182447
$ws.Range("G947").Value = 38684
$ws.Range("H947").Value = 27040
$ws.Range("I947").Value = 66487
$ws.Range("J947").Value = 0
$ws.Range("K947").Value = 0
$ws.Range("L947").Value = 0
$ws.Range("M947").Value = 2
$ws.Range("O947").Value = 2
$ws.Range("P947").Value = 0
$ws.Range("Q947").Value = 2
$ws.Range("T947").Value = 136

# Row 948: 2022/11/01
$cell = $ws.Range("A948")
$cell.Formula = "=TEXT(""2022/11/01"",""@"")"
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Range("B948").Value = 947
$ws.Range("C948").Value = 39448
$ws.Range("D948").Value = 691
$ws.Range("F948").Value = 0.01751673089
$ws.Range("G948").Value = 38685
$ws.Range("H948").Value = 27040
$ws.Range("I948").Value = 66488
$ws.Range("J948").Value = 1
$ws.Range("K948").Value = 1
$ws.Range("L948").Value = 1
$ws.Range("M948").Value = 0
$ws.Range("O948").Value = 0
$ws.Range("P948").Value = 0
$ws.Range("Q948").Value = 2
$ws.Range("T948").Value = 136

$excel.CutCopyMode = $false
Write-Output "Done"